$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
# Row 2
$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:04:21'
# Row 3
$ws.Cells.Item(3, 1).Value = 'Total filas: 228'
# Row 38
$ws.Cells.Item(38, 1).Value = '05:59:00'
$ws.Cells.Item(38, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(38, 4).Value = 66
# Row 39
$ws.Cells.Item(39, 1).Value = '06:50:53'
$ws.Cells.Item(39, 3).Value = '15_ABASTO'
$ws.Cells.Item(39, 4).Value = 15
# Row 53
$ws.Cells.Item(53, 1).Value = '05:59:00'
$ws.Cells.Item(53, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(53, 4).Value = 93
# Row 55
$ws.Cells.Item(55, 1).Value = '06:50:53'
$ws.Cells.Item(55, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(55, 4).Value = 42
# Row 67
$ws.Cells.Item(67, 3).Value = '10_OLMOS'
# Row 68
$ws.Cells.Item(68, 3).Value = '16_SANTA ANA'
# Row 75
$ws.Cells.Item(75, 1).Value = '07:35:06'
$ws.Cells.Item(75, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(75, 4).Value = 48
# Row 76
$ws.Cells.Item(76, 1).Value = '08:22:49'
$ws.Cells.Item(76, 3).Value = '215B_EL PATO'
$ws.Cells.Item(76, 4).Value = 1
# Row 98
$ws.Cells.Item(98, 3).Value = '16_SANTA ANA'
# Row 99
$ws.Cells.Item(99, 3).Value = '17_ROMERO'
# Row 100
$ws.Cells.Item(100, 1).Value = '08:54:22'
$ws.Cells.Item(100, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(100, 4).Value = 29
# Row 101
$ws.Cells.Item(101, 1).Value = '07:35:06'
$ws.Cells.Item(101, 3).Value = '17_ROMERO'
$ws.Cells.Item(101, 4).Value = 108
# Row 106
$ws.Cells.Item(106, 1).Value = '08:54:22'
$ws.Cells.Item(106, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(106, 4).Value = 40
# Row 107
$ws.Cells.Item(107, 1).Value = '08:22:49'
$ws.Cells.Item(107, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(107, 4).Value = 72
# Row 154
$ws.Cells.Item(154, 1).Value = '12:06:53'
$ws.Cells.Item(154, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(154, 4).Value = 0
# Row 156
$ws.Cells.Item(156, 1).Value = '11:27:45'
$ws.Cells.Item(156, 3).Value = '14_ABASTO'
$ws.Cells.Item(156, 4).Value = 39
# Row 166
$ws.Cells.Item(166, 3).Value = '27_EL RETIRO'
# Row 167
$ws.Cells.Item(167, 3).Value = '16_SANTA ANA'
# Row 170
$ws.Cells.Item(170, 1).Value = '12:06:53'
$ws.Cells.Item(170, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(170, 4).Value = 42
# Row 172
$ws.Cells.Item(172, 1).Value = '12:48:55'
$ws.Cells.Item(172, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(172, 4).Value = 0
# Row 198
$ws.Cells.Item(198, 3).Value = '11_ETCHEVERRY'
# Row 199
$ws.Cells.Item(199, 3).Value = '215A_EL PATO'
# Row 201
$ws.Cells.Item(201, 3).Value = '16_P MOR-167 Y 521'
# Row 202
$ws.Cells.Item(202, 3).Value = '225_GOMEZ'
# Row 203
$ws.Cells.Item(203, 3).Value = '23_HERNANDEZ'
# Row 204
$ws.Cells.Item(204, 1).Value = '14:04:21'
$ws.Cells.Item(204, 3).Value = '17_ROMERO'
$ws.Cells.Item(204, 4).Value = 0
# Row 205
$ws.Cells.Item(205, 1).Value = '14:04:21'
$ws.Cells.Item(205, 2).Value = '14:04'
$ws.Cells.Item(205, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(205, 4).Value = 0
# Row 206
$ws.Cells.Item(206, 1).Value = '14:04:21'
$ws.Cells.Item(206, 2).Value = '14:05'
$ws.Cells.Item(206, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(206, 4).Value = 1
# Row 207
$ws.Cells.Item(207, 1).Value = '14:04:21'
$ws.Cells.Item(207, 2).Value = '14:05'
$ws.Cells.Item(207, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(207, 4).Value = 1
# Row 208
$ws.Cells.Item(208, 1).Value = '14:04:21'
$ws.Cells.Item(208, 2).Value = '14:07'
$ws.Cells.Item(208, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(208, 4).Value = 3
# Row 209
$ws.Cells.Item(209, 1).Value = '14:04:21'
$ws.Cells.Item(209, 2).Value = '14:12'
$ws.Cells.Item(209, 3).Value = '15_ABASTO'
$ws.Cells.Item(209, 4).Value = 8
# Row 210
$ws.Cells.Item(210, 1).Value = '13:18:34'
$ws.Cells.Item(210, 2).Value = '14:16'
$ws.Cells.Item(210, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(210, 4).Value = 58
# Row 211
$ws.Cells.Item(211, 1).Value = '14:04:21'
$ws.Cells.Item(211, 2).Value = '14:17'
$ws.Cells.Item(211, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(211, 4).Value = 13
# Row 212
$ws.Cells.Item(212, 1).Value = '12:48:55'
$ws.Cells.Item(212, 2).Value = '14:19'
$ws.Cells.Item(212, 3).Value = '215C_EL PATO'
$ws.Cells.Item(212, 4).Value = 91
# Row 213
$ws.Cells.Item(213, 1).Value = '12:48:55'
$ws.Cells.Item(213, 2).Value = '14:20'
$ws.Cells.Item(213, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(213, 4).Value = 92
# Row 214
$ws.Cells.Item(214, 1).Value = '14:04:21'
$ws.Cells.Item(214, 2).Value = '14:20'
$ws.Cells.Item(214, 3).Value = '215C_EL PATO'
$ws.Cells.Item(214, 4).Value = 16
# Row 215
$ws.Cells.Item(215, 1).Value = '14:04:21'
$ws.Cells.Item(215, 2).Value = '14:21'
$ws.Cells.Item(215, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(215, 4).Value = 17
# Row 216
$ws.Cells.Item(216, 1).Value = '14:04:21'
$ws.Cells.Item(216, 2).Value = '14:44'
$ws.Cells.Item(216, 3).Value = '14_ABASTO'
$ws.Cells.Item(216, 4).Value = 40
# Row 217
$ws.Cells.Item(217, 1).Value = '13:18:34'
$ws.Cells.Item(217, 2).Value = '14:45'
$ws.Cells.Item(217, 3).Value = '14_ABASTO'
$ws.Cells.Item(217, 4).Value = 87
$ws.Cells.Item(217, 5).Value = 'LP1912'
# Row 218
$ws.Cells.Item(218, 1).Value = '14:04:21'
$ws.Cells.Item(218, 2).Value = '14:56'
$ws.Cells.Item(218, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(218, 4).Value = 52
$ws.Cells.Item(218, 5).Value = 'LP1912'
# Row 219
$ws.Cells.Item(219, 1).Value = '14:04:21'
$ws.Cells.Item(219, 2).Value = '14:58'
$ws.Cells.Item(219, 3).Value = '215B_EL PATO'
$ws.Cells.Item(219, 4).Value = 54
$ws.Cells.Item(219, 5).Value = 'LP1912'
# Row 220
$ws.Cells.Item(220, 1).Value = '14:04:21'
$ws.Cells.Item(220, 2).Value = '15:00'
$ws.Cells.Item(220, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(220, 4).Value = 56
$ws.Cells.Item(220, 5).Value = 'LP1912'
# Row 221
$ws.Cells.Item(221, 1).Value = '14:04:21'
$ws.Cells.Item(221, 2).Value = '15:05'
$ws.Cells.Item(221, 3).Value = '10_OLMOS'
$ws.Cells.Item(221, 4).Value = 61
$ws.Cells.Item(221, 5).Value = 'LP1912'
# Row 222
$ws.Cells.Item(222, 1).Value = '14:04:21'
$ws.Cells.Item(222, 2).Value = '15:10'
$ws.Cells.Item(222, 3).Value = '17_ROMERO'
$ws.Cells.Item(222, 4).Value = 66
$ws.Cells.Item(222, 5).Value = 'LP1912'
# Row 223
$ws.Cells.Item(223, 1).Value = '14:04:21'
$ws.Cells.Item(223, 2).Value = '15:13'
$ws.Cells.Item(223, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(223, 4).Value = 69
$ws.Cells.Item(223, 5).Value = 'LP1912'
# Row 224
$ws.Cells.Item(224, 1).Value = '14:04:21'
$ws.Cells.Item(224, 2).Value = '15:20'
$ws.Cells.Item(224, 3).Value = '15_ABASTO'
$ws.Cells.Item(224, 4).Value = 76
$ws.Cells.Item(224, 5).Value = 'LP1912'
# Row 225
$ws.Cells.Item(225, 1).Value = '14:04:21'
$ws.Cells.Item(225, 2).Value = '15:21'
$ws.Cells.Item(225, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(225, 4).Value = 77
$ws.Cells.Item(225, 5).Value = 'LP1912'
# Row 226
$ws.Cells.Item(226, 1).Value = '14:04:21'
$ws.Cells.Item(226, 2).Value = '15:32'
$ws.Cells.Item(226, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(226, 4).Value = 88
$ws.Cells.Item(226, 5).Value = 'LP1912'
# Row 227
$ws.Cells.Item(227, 1).Value = '14:04:21'
$ws.Cells.Item(227, 2).Value = '15:37'
$ws.Cells.Item(227, 3).Value = '10_OLMOS'
$ws.Cells.Item(227, 4).Value = 93
$ws.Cells.Item(227, 5).Value = 'LP1912'
# Row 228
$ws.Cells.Item(228, 1).Value = '14:04:21'
$ws.Cells.Item(228, 2).Value = '15:38'
$ws.Cells.Item(228, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(228, 4).Value = 94
$ws.Cells.Item(228, 5).Value = 'LP1912'
# Row 229
$ws.Cells.Item(229, 1).Value = '14:04:21'
$ws.Cells.Item(229, 2).Value = '15:38'
$ws.Cells.Item(229, 3).Value = '215A_EL PATO'
$ws.Cells.Item(229, 4).Value = 94
$ws.Cells.Item(229, 5).Value = 'LP1912'
# Row 230
$ws.Cells.Item(230, 1).Value = '14:04:21'
$ws.Cells.Item(230, 2).Value = '15:46'
$ws.Cells.Item(230, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(230, 4).Value = 102
$ws.Cells.Item(230, 5).Value = 'LP1912'
# Row 231
$ws.Cells.Item(231, 1).Value = '14:04:21'
$ws.Cells.Item(231, 2).Value = '15:50'
$ws.Cells.Item(231, 3).Value = '14_ABASTO'
$ws.Cells.Item(231, 4).Value = 106
$ws.Cells.Item(231, 5).Value = 'LP1912'
# Row 232
$ws.Cells.Item(232, 1).Value = '14:04:21'
$ws.Cells.Item(232, 2).Value = '15:54'
$ws.Cells.Item(232, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(232, 4).Value = 110
$ws.Cells.Item(232, 5).Value = 'LP1912'
# Row 233
$ws.Cells.Item(233, 1).Value = '14:04:21'
$ws.Cells.Item(233, 2).Value = '15:57'
$ws.Cells.Item(233, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(233, 4).Value = 113
$ws.Cells.Item(233, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
# Row 2
$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:04:21'
# Row 3
$ws.Cells.Item(3, 1).Value = 'Total filas: 28'
# Row 31
$ws.Cells.Item(31, 1).Value = '14:04:21'
$ws.Cells.Item(31, 4).Value = 16
# Row 32
$ws.Cells.Item(32, 1).Value = '14:04:21'
$ws.Cells.Item(32, 4).Value = 54
# Row 33
$ws.Cells.Item(33, 1).Value = '14:04:21'
$ws.Cells.Item(33, 2).Value = '15:38'
$ws.Cells.Item(33, 3).Value = '215A_EL PATO'
$ws.Cells.Item(33, 4).Value = 94
$ws.Cells.Item(33, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
# Row 2
$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:04:21'
# Row 3
$ws.Cells.Item(3, 1).Value = 'Total filas: 32'
# Row 35
$ws.Cells.Item(35, 1).Value = '14:04:21'
$ws.Cells.Item(35, 2).Value = '14:11'
$ws.Cells.Item(35, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(35, 4).Value = 7
$ws.Cells.Item(35, 5).Value = 'L6173'
# Row 36
$ws.Cells.Item(36, 1).Value = '14:04:21'
$ws.Cells.Item(36, 2).Value = '14:53'
$ws.Cells.Item(36, 3).Value = '215D_LA PLATA'
$ws.Cells.Item(36, 4).Value = 49
$ws.Cells.Item(36, 5).Value = 'L6203'
# Row 37
$ws.Cells.Item(37, 1).Value = '14:04:21'
$ws.Cells.Item(37, 2).Value = '15:34'
$ws.Cells.Item(37, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(37, 4).Value = 90
$ws.Cells.Item(37, 5).Value = 'L6173'
